$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# optimization_parameters sheet (7th sheet, originally tab #7 in the book)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1 had duplicated "value" header cells in C1:F1 - remove them so the
# header row only spans A1:B1.
$ws.Range("C1:F1").ClearContents()

# "Model" label is renamed to "production_function" (value "Sigmoid" stays).
$ws.Range("A8").Value = "production_function"

# Insert a brand new "L_curve" parameter row right after "production_function",
# pushing everything below it down by one row.
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0
$ws.Range("B9").NumberFormat = "0.00E+00"

# The old "Strain"/"Deletion" mapping row (now at row 17 after the insert
# above) is removed entirely.
$ws.Rows.Item(17).Delete()

# Make this the active sheet/tab, with C1:F1 selected (matches new layout).
$ws.Activate()
$ws.Range("C1:F1").Select()
